# Premier League - atualizacao de dados
# Atualizando rodada 5 do sabado - 22 de setembro
# Colunas B..P: partidas jogadas, partidas mandante, partidas visitante,
# gols mandante, media gols mandante, gols visitante, media gols visitante,
# total gols, media gols, gols sofridos mandante, media gols sofridos mand,
# gols sofridos visitante, media gols sofridos visit, total gols sofridos, media gols sofridos

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{}
$rowData[3] = @(5, 3, 2, 6, 1.2, 4, 2, 10, 2, 5, 1.666666666666667, 2, 1, 7, 1.4)  # Aston Villa
$rowData[4] = @(5, 2, 3, 1, 0.2, 4, 1.333333333333333, 5, 1, 2, 1, 6, 2, 8, 1.6)  # Bournemouth
$rowData[5] = @(5, 2, 3, 5, 1, 2, 0.6666666666666666, 7, 1.4, 2, 1, 7, 2.333333333333333, 9, 1.8)  # Brentford
$rowData[7] = @(5, 2, 3, 1, 0.2, 10, 3.333333333333333, 11, 2.2, 3, 1.5, 2, 0.6666666666666666, 5, 1)  # Chelsea
$rowData[9] = @(5, 3, 2, 2, 0.4, 2, 1, 4, 0.8, 4, 1.333333333333333, 3, 1.5, 7, 1.4)  # Crystal Palace
$rowData[10] = @(5, 2, 3, 2, 0.4, 3, 1, 5, 1, 6, 3, 8, 2.666666666666667, 14, 2.8)  # Everton
$rowData[11] = @(5, 3, 2, 6, 1.2, 1, 0.5, 7, 1.4, 3, 1, 2, 1, 5, 1)  # Fulham
$rowData[12] = @(5, 2, 3, 1, 0.2, 2, 0.6666666666666666, 3, 0.6, 3, 1.5, 5, 1.666666666666667, 8, 1.6)  # Ipswich Town
$rowData[13] = @(5, 3, 2, 3, 0.6, 3, 1.5, 6, 1.2, 4, 1.333333333333333, 4, 2, 8, 1.6)  # Leicester City
$rowData[14] = @(5, 3, 2, 5, 1, 5, 2.5, 10, 2, 1, 0.3333333333333333, 0, 0, 1, 0.2)  # Liverpool
$rowData[15] = @(5, 2, 3, 1, 0.2, 4, 1.333333333333333, 5, 1, 3, 1.5, 2, 0.6666666666666666, 5, 1)  # Manchester United
$rowData[16] = @(5, 2, 3, 3, 0.6, 4, 1.333333333333333, 7, 1.4, 1, 0.5, 5, 1.666666666666667, 6, 1.2)  # Newcastle
$rowData[18] = @(5, 3, 2, 1, 0.2, 1, 0.5, 2, 0.4, 5, 1.666666666666667, 4, 2, 9, 1.8)  # Southampton
$rowData[19] = @(5, 3, 2, 7, 1.4, 2, 1, 9, 1.8, 2, 0.6666666666666666, 3, 1.5, 5, 1)  # Tottenham
$rowData[20] = @(5, 3, 2, 2, 0.4, 3, 1.5, 5, 1, 8, 2.666666666666667, 1, 0.5, 9, 1.8)  # West Ham
$rowData[21] = @(5, 2, 3, 3, 0.6, 2, 0.6666666666666666, 5, 1, 4, 2, 6, 2, 14, 2.8)  # Wolves

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}

Write-Host "Updated round 5 data for 16 clubs"
